$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Config Name changes from "Test" to "Debug" (values unchanged)
$ws.Range("A2").Value = "Debug"

# Row 3: Config Name changes from "MyConfig" to "Fast", plus several numeric updates
$ws.Range("A3").Value = "Fast"
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 10
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = 10
$ws.Range("G3").Value = 40
$ws.Range("J3").Value = 3

# Update the active selection to match the new cursor position
$ws.Range("H8").Select()
